$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.332.66'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '3.077.20'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").Value = "'171.37"
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.073.88'
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("D10").Value = "'6.26"
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("E11").Value = '  -2.14%  '
$ws.Range("E12").Value = '  -2.33%  '
$ws.Range("E13").Value = '  -3.57%  '
$ws.Range("D14").Value = "'35.85"
$ws.Range("E14").Value = '  -3.44%  '
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").Value = '3.590.57'
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("D17").Value = '66.325.09'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("E18").Value = '  -2.68%  '
$ws.Range("D19").Value = "'16.59"
$ws.Range("E19").Value = '  +2.05%  '
$ws.Range("D20").Value = '3.077.02'
$ws.Range("D21").Value = "'487.03"
$ws.Range("E21").Value = '  +2.48%  '
$ws.Range("E22").Value = '  -3.16%  '
$ws.Range("E23").Value = '  -2.61%  '
$ws.Range("D24").Value = "'82.38"
$ws.Range("E24").Value = '  -1.70%  '
$ws.Range("D25").Value = "'12.61"
$ws.Range("E25").Value = '  -4.73%  '
$ws.Range("E26").Value = '  -3.08%  '
$ws.Range("D27").Value = "'10.09"
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").Value = "'7.86"
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").Value = "'2.24"
$ws.Range("E30").Value = '  -5.23%  '
$ws.Range("D31").Value = "'2.59"
$ws.Range("E31").Value = '  -3.11%  '
$ws.Range("E32").Value = '  -3.17%  '
$ws.Range("D33").Value = "'0.112"
$ws.Range("E33").Value = '  -2.71%  '
$ws.Range("E34").Value = '  -3.65%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = "'47.89"
$ws.Range("E36").Value = '  +2.21%  '
$ws.Range("E37").Value = '  -4.68%  '
$ws.Range("D38").Value = "'0.942"
$ws.Range("E38").Value = '  -3.47%  '
$ws.Range("D39").Value = "'0.122"
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("D40").Value = "'0.302"
$ws.Range("E40").Value = '  -3.06%  '
$ws.Range("E41").Value = '  -4.66%  '
$ws.Range("D42").Value = "'8.22"
$ws.Range("E42").Value = '  -4.44%  '
$ws.Range("D43").Value = '2.769.49'
$ws.Range("E43").Value = '  -1.68%  '
$ws.Range("D44").Value = "'2.54"
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("E45").Value = '  -2.49%  '
$ws.Range("D46").Value = "'134.76"
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").Value = "'364.87"
$ws.Range("E47").Value = '  -4.49%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  -2.78%  '
$ws.Range("E50").Value = '  -2.41%  '
$ws.Range("E51").Value = '  -2.14%  '
